$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last refreshed" timestamp in the title cell (A1)
$ws.Range("A1").Value2 = 'Datos actualizados a 3 de Julio de 2020 a las 23:40'

# Refreshed COVID-19 snapshot: new totals pulled in for several countries.
# A handful of rows (Barein, Congo, Ruanda, Dominica, Groenlandia, ...) also
# overtook their neighbours in "Casos totales" (or, for two same-count ties,
# simply swapped order), so the table's ranking shifts for those rows too.
# Each tuple below is: row, Pais, Casos totales, Nuevos casos, Casos activos,
# Recuperados, Casos criticos, Muertes hoy, Muertes - i.e. columns A-H.
$updates = @(
    @(4, 'Estados Unidos', 2879926, 44242, 1204693, 1543256, 0, 492, 131977),
    @(5, 'Brasil', 1539081, 37728, 916147, 559760, 0, 1184, 63174),
    @(9, 'Peru', 295599, 3595, 185852, 99521, 0, 181, 10226),
    @(20, 'Francia', 166960, 582, 77060, 60007, 0, 18, 29893),
    @(49, 'Barein', 28410, 573, 23318, 4997, 0, 1, 95),
    @(50, 'Rumania', 28166, 420, 19545, 6913, 0, 21, 1708),
    @(51, 'Israel', 28055, 1008, 17669, 10060, 0, 2, 326),
    @(122, 'Congo', 1557, 175, 501, 1012, 0, 3, 44),
    @(123, 'Nueva Zelanda', 1530, 0, 1490, 18, 0, 0, 22),
    @(124, 'Sierra Leona', 1524, 6, 1042, 420, 0, 2, 62),
    @(125, 'Malaui', 1402, 60, 317, 1069, 0, 0, 16),
    @(133, 'Ruanda', 1081, 18, 512, 566, 0, 0, 3),
    @(134, 'Niger', 1081, 0, 959, 54, 0, 0, 68),
    @(146, 'Santo Tome y Principe', 719, 2, 267, 439, 0, 0, 13),
    @(154, 'Surinam', 561, 14, 267, 281, 0, 0, 13),
    @(166, 'Guyana', 256, 6, 117, 125, 0, 0, 14),
    @(167, 'Martinica', 249, 7, 98, 137, 0, 0, 14),
    @(173, 'Guadalupe', 184, 2, 157, 13, 0, 0, 14),
    @(205, 'Dominica', 18, 0, 18, 0, 0, 0, 0),
    @(206, 'Fiyi', 18, 0, 18, 0, 0, 0, 0),
    @(209, 'Groenlandia', 13, 0, 13, 0, 0, 0, 0),
    @(210, 'Islas Malvinas', 13, 0, 13, 0, 0, 0, 0)
)

foreach ($u in $updates) {
    $r = $u[0]
    for ($c = 1; $c -le 8; $c++) {
        $ws.Cells.Item($r, $c).Value2 = $u[$c]
    }
}

Write-Output "Updated $($updates.Count) data rows plus title."
